{"js": "// Remove the unused \"Abstract Title\" paragraph style and the unused\n// \"Footnote Block Text\" paragraph style, and bump the \"Abstract\" style's\n// space-before from 5pt (100 twips) to 15pt (300 twips) so it matches\n// the space-after value.\n\n// Delete the \"Abstract Title\" style entirely.\nconst abstractTitleStyle = context.document.getStyles().getByNameOrNullObject(\"Abstract Title\");\nawait context.sync();\nabstractTitleStyle.delete();\nawait context.sync();\n\n// Delete the \"Footnote Block Text\" style entirely.\nconst footnoteBlockTextStyle = context.document.getStyles().getByNameOrNullObject(\"Footnote Block Text\");\nawait context.sync();\nfootnoteBlockTextStyle.delete();\nawait context.sync();\n\n// \"Abstract\" style: change paragraph spacing before from 5pt to 15pt (100 -> 300 twentieths of a point).\nconst abstractStyle = context.document.getStyles().getByNameOrNullObject(\"Abstract\");\nawait context.sync();\nabstractStyle.paragraphFormat.spaceBefore = 15;\nawait context.sync();\n", "ps1": "# Remove the unused \"Abstract Title\" paragraph style and the unused\n# \"Footnote Block Text\" paragraph style, and bump the \"Abstract\" style's\n# space-before from 5pt (100 twips) to 15pt (300 twips) so it matches\n# the space-after value.\n\n$d = $word.ActiveDocument\n\n# Delete the \"Abstract Title\" style entirely.\n$d.Styles(\"Abstract Title\").Delete()\n\n# \"Abstract\" style: change paragraph spacing before from 100 (5pt) to 300 (15pt) twentieths of a point.\n$d.Styles(\"Abstract\").ParagraphFormat.SpaceBefore = 15\n\n# Delete the \"Footnote Block Text\" style entirely.\n$d.Styles(\"Footnote Block Text\").Delete()\n"}
